# Update the "想去人数" (want-to-go count) figures (column F) that changed
# between scrape runs, across the four worksheets of the workbook.
# Mapping derived from the OOXML diff: row -> (old value, new value)

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$updates1 = @{
    4  = 313
    5  = 1266
    7  = 302
    8  = 1103
    9  = 429
    10 = 6909
    13 = 30
    14 = 7809
    17 = 5267
    19 = 2288
    20 = 971
    22 = 255
    26 = 294
    27 = 231
    29 = 2026
    31 = 228
    33 = 536
    34 = 22
    35 = 1381
    37 = 2112
    39 = 14
}
foreach ($row in $updates1.Keys) {
    $ws1.Range("F$row").Value = $updates1[$row]
}

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value = 90

# Sheet "本地生活"
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 1260

# Sheet "全部类型" (combined listing of all the above)
$ws4 = $wb.Worksheets.Item("全部类型")
$updates4 = @{
    4  = 1260
    7  = 313
    8  = 1266
    11 = 302
    12 = 1103
    13 = 429
    14 = 6909
    17 = 30
    18 = 7809
    21 = 5267
    23 = 2288
    24 = 971
    26 = 255
    32 = 294
    33 = 231
    35 = 2026
    37 = 228
    39 = 536
    40 = 22
    42 = 1381
    44 = 2112
    47 = 14
    48 = 90
}
foreach ($row in $updates4.Keys) {
    $ws4.Range("F$row").Value = $updates4[$row]
}
